# Add a new "Senior plus plus" tier row (row 25) to the DiscountRules
# decision table, following the same A=name / B:E=repeated-value pattern
# used by the existing rows (e.g. row 24 "Senior plus" / 45).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A25").Value = "Senior plus plus"

# Force the value cells to be stored as text ("50"), matching the
# original sheet where every cell (including numeric-looking ones) uses
# shared-string "t=s" entries rather than numeric values. Temporarily
# apply a text number format so Excel doesn't auto-convert "50" to a
# number, then restore the default "Normal" style so the new cells line
# up with the rest of the sheet (which uses the default style).
$valueCells = $ws.Range("B25:E25")
$valueCells.NumberFormat = "@"
$valueCells.Value = "50"
$valueCells.Style = "Normal"
